$d = $word.ActiveDocument
$ps = $d.PageSetup
$ps.Orientation = 0
$ps.BottomMargin = 144.0
